$d = $word.ActiveDocument

# Make sure the anchor paragraph is really there before we touch anything.
$check = $d.Content
$check.Find.ClearFormatting()
$found = $check.Find.Execute("A ni bil ta tip učas bl šlank?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph 'A ni bil ta tip učas bl šlank?' was not found."
}

# Locate the (last) paragraph whose text is "A ni bil ta tip učas bl šlank?"
# so we can insert the new paragraph right after it.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*A ni bil ta tip učas bl šlank?*") {
        $targetPara = $p
    }
}
if ($null -eq $targetPara) {
    throw "Could not re-locate the anchor paragraph via Paragraphs collection."
}

# Add a brand-new paragraph right after it.
$tr = $targetPara.Range
$tr.Collapse(0)
$tr.InsertParagraphAfter()

$newPara = $targetPara.Next()
$nr = $newPara.Range
$nr.Collapse(1)
$newParaStart = $nr.Start

# Type the whole sentence in one go so it correctly inherits the run
# formatting (the en-US language mark) from the surrounding paragraph.
$nr.InsertAfter("Aja, sej to ni ta.")

# Split "Aja, sej to ni ta." into three runs -- "Aja, sej to ni" / " t" / "a." --
# by nudging (and then reverting) character formatting on the middle span,
# which forces Word to break it into separate runs without altering the
# final rendered formatting of any of them.
$splitStart = $newParaStart + "Aja, sej to ni".Length
$splitEnd = $splitStart + " t".Length
$mid = $d.Range($splitStart, $splitEnd)
$mid.Bold = 1
$mid.Bold = 0
